$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of date values (Excel serial 42036 == 2015-02-01).
$ws.Range("A5:C5").Value = 42036

# Apply the built-in date number format to A5, then copy that
# formatting (and only the formatting) across the rest of the row so
# every cell shares a single reused cell-style (xf) entry.
$ws.Range("A5").NumberFormat = "mm-dd-yy"
$ws.Range("A5").Copy()
$ws.Range("B5:C5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the selection to C5 to match the saved cursor position.
$ws.Range("C5").Select()
